$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2790.8462
$ws.Range("I132").Value = 2595.0715
$ws.Range("J132").Value = 3019.25
$ws.Range("K132").Value = 7785.2145
$ws.Range("L132").Value = 9057.75
$ws.Range("M132").Value = -5255.2145
$ws.Range("N132").Value = -14117.75
$ws.Range("H137").Value = 30715.068
$ws.Range("I137").Value = 36510.832
$ws.Range("J137").Value = 2895.4
$ws.Range("K137").Value = 109532.496
$ws.Range("L137").Value = 8686.200000000001
$ws.Range("M137").Value = -106982.496
$ws.Range("N137").Value = -13786.2
$ws.Range("H138").Value = 2549.63
$ws.Range("J138").Value = 2681.4888
$ws.Range("L138").Value = 8044.4664
$ws.Range("N138").Value = -18324.4664
$ws.Range("H141").Value = 7047.722
$ws.Range("I141").Value = 5868.125
$ws.Range("J141").Value = 9406.916999999999
$ws.Range("K141").Value = 17604.375
$ws.Range("L141").Value = 28220.751
$ws.Range("M141").Value = -12424.375
$ws.Range("N141").Value = -38580.751

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 116000
$ws.Range("J122").Value = 116000
$ws.Range("L122").Value = 116000
$ws.Range("N122").Value = -125800
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 99999.5
$ws.Range("J126").Value = 99999.5
$ws.Range("L126").Value = 99999.5
$ws.Range("N126").Value = -109879.5
$ws.Range("H129").Value = 113997.5
$ws.Range("J129").Value = 113997.5
$ws.Range("L129").Value = 113997.5
$ws.Range("N129").Value = -123997.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 540.8333
$ws.Range("I22").Value = 311.5
$ws.Range("K22").Value = 311.5
$ws.Range("M22").Value = 38.5
$ws.Range("H31").Value = 8797.423000000001
$ws.Range("I31").Value = 4599.7
$ws.Range("J31").Value = 11421
$ws.Range("K31").Value = 4599.7
$ws.Range("L31").Value = 11421
$ws.Range("M31").Value = -4304.7
$ws.Range("N31").Value = -12011
$ws.Range("H34").Value = 8797.423000000001
$ws.Range("I34").Value = 4599.7
$ws.Range("J34").Value = 11421
$ws.Range("K34").Value = 4599.7
$ws.Range("L34").Value = 11421
$ws.Range("M34").Value = -4397.7
$ws.Range("N34").Value = -11825

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 800549.3
$ws.Range("I11").Value = 1000285
$ws.Range("K11").Value = 3000855
$ws.Range("M11").Value = -3000715
$ws.Range("H16").Value = 599.6667
$ws.Range("I16").Value = 699.5
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 2098.5
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -1925.5
$ws.Range("N16").Value = -1546
$ws.Range("H18").Value = 3174.5454
$ws.Range("I18").Value = 987
$ws.Range("K18").Value = 2961
$ws.Range("M18").Value = -2792
$ws.Range("H26").Value = 170
$ws.Range("I26").Value = 162.5
$ws.Range("K26").Value = 487.5
$ws.Range("M26").Value = -199.5
$ws.Range("N26").Value = -1176
$ws.Range("H34").Value = 1043.625
$ws.Range("I34").Value = 86.25
$ws.Range("J34").Value = 2001
$ws.Range("K34").Value = 258.75
$ws.Range("L34").Value = 6003
$ws.Range("M34").Value = -174.75
$ws.Range("N34").Value = -6171
$ws.Range("H39").Value = 3162.2
$ws.Range("J39").Value = 3752.75
$ws.Range("L39").Value = 11258.25
$ws.Range("N39").Value = -11846.25
$ws.Range("H55").Value = 3324.3333
$ws.Range("I55").Value = 1224.75
$ws.Range("J55").Value = 5004
$ws.Range("K55").Value = 3674.25
$ws.Range("L55").Value = 15012
$ws.Range("M55").Value = -3497.25
$ws.Range("N55").Value = -15366
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H103").Value = 849
$ws.Range("I103").Value = 849
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 2547
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -1668
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 1199.6666
$ws.Range("I132").Value = 697
$ws.Range("J132").Value = 1601.8
$ws.Range("K132").Value = 6273
$ws.Range("L132").Value = 14416.2
$ws.Range("M132").Value = -3743
$ws.Range("N132").Value = -19476.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1425.2222
$ws.Range("J122").Value = 1799.75
$ws.Range("L122").Value = 5399.25
$ws.Range("N122").Value = -10299.25
$ws.Range("H126").Value = 2676.6667
$ws.Range("I126").Value = 2671.818
$ws.Range("J126").Value = 2682
$ws.Range("K126").Value = 8015.454000000001
$ws.Range("L126").Value = 8046
$ws.Range("M126").Value = -5545.454000000001
$ws.Range("N126").Value = -12986

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2558.9167
$ws.Range("I22").Value = 1315
$ws.Range("K22").Value = 1315
$ws.Range("M22").Value = -1020
$ws.Range("H27").Value = 2558.9167
$ws.Range("I27").Value = 1315
$ws.Range("K27").Value = 1315
$ws.Range("M27").Value = -1208
$ws.Range("H46").Value = 8974.046
$ws.Range("I46").Value = 4700
$ws.Range("J46").Value = 9177.571
$ws.Range("K46").Value = 4700
$ws.Range("L46").Value = 9177.571
$ws.Range("M46").Value = -4512
$ws.Range("N46").Value = -9553.571
$ws.Range("H122").Value = 12793.526
$ws.Range("I122").Value = 13415.777
$ws.Range("K122").Value = 40247.331
$ws.Range("M122").Value = -37797.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1916.5834
$ws.Range("I126").Value = 1299.9
$ws.Range("K126").Value = 3899.7
$ws.Range("M126").Value = -1429.7
$ws.Range("H132").Value = 2250.1875
$ws.Range("I132").Value = 2178.8572
$ws.Range("J132").Value = 2749.5
$ws.Range("K132").Value = 6536.571599999999
$ws.Range("L132").Value = 8248.5
$ws.Range("M132").Value = -4006.571599999999
$ws.Range("N132").Value = -13308.5
